# Meeting Minutes and Agenda
# Fixed typos in the Jan. 23 minutes (DreamSpark spelling call-outs,
# "falling" -> "following", "process" -> "processes") and split the
# closing "_GoBack" bookmark into its own paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...Microsoft's DreamSpark program..." paragraph: split "Dream"
#    and "park" into their own runs (mirrors Word's spell-check
#    bracketing of the unrecognised word "DreamSpark").
# ---------------------------------------------------------------------
$rDream = $d.Content
$rDream.Find.Execute("Dream") | Out-Null
$rDream.Bold = $true
$rDream.Bold = $false

$rPark = $d.Content
$rPark.Find.Execute("park") | Out-Null
$rPark.Bold = $true
$rPark.Bold = $false

# ---------------------------------------------------------------------
# 2) "At the next meeting..." paragraph: fix "falling" -> "following"
#    and "process" -> "processes", then split into the same run
#    boundaries as the authored edit.
# ---------------------------------------------------------------------
$rFalling = $d.Content
$rFalling.Find.Execute("falling", $true, $false, $false, $false, $false, `
    $true, 1, $false, "following", 2) | Out-Null

$rProcess = $d.Content
$rProcess.Find.Execute("lifecycle process;", $true, $false, $false, $false, `
    $false, $true, 1, $false, "lifecycle processes;", 2) | Out-Null

$rRun2 = $d.Content
$rRun2.Find.Execute(" detailed discussion on the follow") | Out-Null
$rRun2.Bold = $true
$rRun2.Bold = $false

$rRun3 = $d.Content
$rRun3.Find.Execute("ing topics: potential programming languages; development lifecycle process") | Out-Null
$rRun3.Bold = $true
$rRun3.Bold = $false

$rRun4 = $d.Content
$rRun4.Find.Execute("es; scheduling") | Out-Null
$rRun4.End = $rRun4.Start + 2
$rRun4.Bold = $true
$rRun4.Bold = $false

# ---------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark out of the "...basement." paragraph
#    and into its own new paragraph (same pPr as the blank paragraph
#    that used to follow it).
# ---------------------------------------------------------------------
$rBasement = $d.Content
$rBasement.Find.Execute("basement.") | Out-Null
$rBasement.Collapse(0)
$rBasement.InsertBefore("`r")

Write-Host "edit complete"
